$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (re-loaded meter_readings table)
$data = @(
    @("01816777", "014673", "kWh", 44848.428611111114),
    @("01809383", "039672", "kWh", 44848.427488425928),
    @("N-002465", "021043", "kWh", 44838.430069444446),
    @("01816775", "018903", "kWh", 44838.429224537038),
    @("01810141", "011511", "kWh", 44838.433668981481),
    @("01810131", "013139", "kWh", 44838.43440972222),
    @("01810125", "022939", "kWh", 44838.433298611111),
    @("01815060", "016060", "kWh", 44838.433506944442),
    @("01810138", "031135", "kWh", 44838.433993055558),
    @("01816787", "023521", "kWh", 44838.434247685182),
    @("016067",   "090965", "kWh", 44838.434618055559),
    @("01810142", "018289", "kWh", 44838.924074074072),
    @("01809380", "004420", "kWH", 44838.428194444445)
)

# Update header label for column C
$ws.Range("C1").Value = "reading_unit"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

$ws.Range("C1").Select()
